$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Overview" sheet: row 3 (b.md) is now ready for handoff.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 02:38:09"

# ---------------------------------------------------------------------
# "zh-cn" sheet: row 3 (b.md) has a fresh handoff package generated.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 02:37:59"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bba672c47f80aacd22804baf2c9708a431a43c3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71095593ffd06658f079307103373f062b3a5efa/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# "de-de" sheet: row 3 (b.md) has a fresh handoff package generated.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 02:38:09"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bba672c47f80aacd22804baf2c9708a431a43c3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71095593ffd06658f079307103373f062b3a5efa/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
